# Update the cumulative best-fitness values in column C (rows 2-252).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10 get individually updated values
$ws.Range("C2").Value  = 12376
$ws.Range("C3").Value  = 12376
$ws.Range("C4").Value  = 12376
$ws.Range("C5").Value  = 12115
$ws.Range("C6").Value  = 11917
$ws.Range("C7").Value  = 11447
$ws.Range("C8").Value  = 11447
$ws.Range("C9").Value  = 10930
$ws.Range("C10").Value = 10930

# Rows 11-252 all flatten out to the same value
$ws.Range("C11:C252").Value = 10453
